# Applies the update described by the commit:
#  - Inserts a new row (row 9) above the old row 9, shifting the rows
#    that followed (old rows 9-14) down to rows 10-15.
#  - The new row 9 gets values: A=5273105816, B="63CFMAR1BR", C=3, D="ZW07"
#  - Cell C8 changes from 5 to 2
#  - Selection moves to A9

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Update C8: 5 -> 2
$ws.Range("C8").Value = 2

# Insert a new row above current row 9 (shifts rows 9:14 down to 10:15)
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with its data
$ws.Cells.Item(9, 1).Value = 5273105816
$ws.Cells.Item(9, 2).Value = "63CFMAR1BR"
$ws.Cells.Item(9, 3).Value = 3
$ws.Cells.Item(9, 4).Value = "ZW07"

# Update the selection to match the saved view state (A9)
$ws.Range("A9").Select()
